$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.217.91"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "'2.994.43"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'501.98"
$ws.Range("E5").Value = "  -4.35%  "
$ws.Range("D6").Value = "'138.83"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.432"
$ws.Range("E8").Value = "  -3.20%  "
$ws.Range("D9").Value = "'7.31"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'0.108"
$ws.Range("E10").Value = "  -4.02%  "
$ws.Range("D11").Value = "'0.361"
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").Value = "'3.507.46"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("D14").Value = "'26.29"
$ws.Range("E14").Value = "  -1.95%  "
$ws.Range("D15").Value = "'0.0000161"
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").Value = "'57.267.48"
$ws.Range("D17").Value = "'6.10"
$ws.Range("E17").Value = "  -2.69%  "
$ws.Range("D18").Value = "'2.997.58"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").Value = "'7.90"
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("D21").Value = "'321.61"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'5.73"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "'0.494"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'63.57"
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  -5.36%  "
$ws.Range("D28").Value = "'0.0₃0901"
$ws.Range("E28").Value = "  -7.23%  "
$ws.Range("D29").Value = "'6.57"
$ws.Range("E29").Value = "  -5.86%  "
$ws.Range("D30").Value = "'7.16"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("E31").Value = "  -4.03%  "
$ws.Range("D32").Value = "'1.17"
$ws.Range("E32").Value = "  -5.62%  "
$ws.Range("D33").Value = "'20.33"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").Value = "'155.16"
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("D35").Value = "'4.59"
$ws.Range("E35").Value = "  -3.54%  "
$ws.Range("D36").Value = "'5.80"
$ws.Range("E36").Value = "  -2.21%  "
$ws.Range("E37").Value = "  -5.92%  "
$ws.Range("D38").Value = "'24.40"
$ws.Range("E38").Value = "  -6.07%  "
$ws.Range("D39").Value = "'0.0666"
$ws.Range("E39").Value = "  -5.64%  "
$ws.Range("B40").Value = "'RenzoRestakedETH"
$ws.Range("C40").Value = "'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "'3.026.42"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("B41").Value = "'OKB"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'37.82"
$ws.Range("E41").Value = "  +0.28%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "'3.77"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("D44").Value = "'0.645"
$ws.Range("E44").Value = "  -2.36%  "
$ws.Range("B45").Value = "'Maker"
$ws.Range("C45").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "'2.207.38"
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("B46").Value = "'Stacks"
$ws.Range("C46").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.40"
$ws.Range("E46").Value = "  -5.76%  "
$ws.Range("D47").Value = "'5.98"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").Value = "'0.945"
$ws.Range("E48").Value = "  -8.36%  "
$ws.Range("D49").Value = "'0.0236"
$ws.Range("E49").Value = "  -4.66%  "
$ws.Range("D50").Value = "'19.38"
$ws.Range("E50").Value = "  -4.39%  "
$ws.Range("D51").Value = "'1.80"
$ws.Range("E51").Value = "  -11.13%  "
